$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.273556709289551
$ws.Range("B1").Value = 3.308527231216431
$ws.Range("C1").Value = 2.905781030654907
$ws.Range("D1").Value = 3.156487464904785
$ws.Range("E1").Value = 2.299679279327393
